$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- A7: update the "as of" date in the confidentiality disclosure text ---
# The sheet is protected and A7 is a locked cell, so briefly unlock it,
# write the new text, then restore its original formatting (no explicit
# style, default row height) and re-lock it so the worksheet protection
# and formatting are left exactly as they were.
$ws.Range("A7").Locked = $false
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-31 for illustrative purposes only and are subject to change."
$ws.Rows.Item(7).AutoFit()
$ws.Range("A8").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# --- D2/E2, D3/E3, E4: updated weight / percent-change figures ---
# Same unlock -> write -> restore-format pattern, using D4 (untouched,
# same original style) as the formatting donor for each changed cell.
$ws.Range("D2").Locked = $false
$ws.Range("E2").Locked = $false
$ws.Range("D3").Locked = $false
$ws.Range("E3").Locked = $false
$ws.Range("E4").Locked = $false

$ws.Range("D2").Value = 0.8433681600850009
$ws.Range("E2").Value = -0.00433070866141716
$ws.Range("D3").Value = 0.1566318399149991
$ws.Range("E3").Value = 0.005087620124364145
$ws.Range("E4").Value = -0.002855498494775732

$ws.Range("D4").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
